$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cck"
$ws.Range("C2").Value = "Cckar"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1295
$ws.Range("H2").Value = 0.3885
$ws.Range("I2").Value = 0.03378978388280364
$ws.Range("J2").Value = 0.03378978388280364
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4642756666666667
$ws.Range("N2").Value = 1.392827
$ws.Range("O2").Value = 0.5486456655941866
$ws.Range("P2").Value = 0.5486456655941866
$ws.Range("Q2").Value = 0.06012369883333334
$ws.Range("R2").Value = 0.5411132895
$ws.Range("S2").Value = 0.01853861846866452
$ws.Range("T2").Value = 0.01853861846866452

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cck"
$ws.Range("C3").Value = "Cckar"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1295
$ws.Range("H3").Value = 0.3885
$ws.Range("I3").Value = 0.03378978388280364
$ws.Range("J3").Value = 0.03378978388280364
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3819456666666667
$ws.Range("N3").Value = 1.145837
$ws.Range("O3").Value = 0.4513543344058135
$ws.Range("P3").Value = 0.4513543344058135
$ws.Range("Q3").Value = 0.04946196383333334
$ws.Range("R3").Value = 0.4451576745
$ws.Range("S3").Value = 0.01525116541413912
$ws.Range("T3").Value = 0.01525116541413912

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cck"
$ws.Range("C4").Value = "Cckar"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.300224333333333
$ws.Range("H4").Value = 3.900673
$ws.Range("I4").Value = 0.3392609978571102
$ws.Range("J4").Value = 0.3392609978571102
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4642756666666667
$ws.Range("N4").Value = 1.392827
$ws.Range("O4").Value = 0.5486456655941866
$ws.Range("P4").Value = 0.5486456655941866
$ws.Range("Q4").Value = 0.6036625191745555
$ws.Range("R4").Value = 5.432962672571
$ws.Range("S4").Value = 0.1861340759794621
$ws.Range("T4").Value = 0.1861340759794621

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cck"
$ws.Range("C5").Value = "Cckar"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.300224333333333
$ws.Range("H5").Value = 3.900673
$ws.Range("I5").Value = 0.3392609978571102
$ws.Range("J5").Value = 0.3392609978571102
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3819456666666667
$ws.Range("N5").Value = 1.145837
$ws.Range("O5").Value = 0.4513543344058135
$ws.Range("P5").Value = 0.4513543344058135
$ws.Range("Q5").Value = 0.4966150498112222
$ws.Range("R5").Value = 4.469535448300999
$ws.Range("S5").Value = 0.1531269218776481
$ws.Range("T5").Value = 0.1531269218776481

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cck"
$ws.Range("C6").Value = "Cckar"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1462553333333333
$ws.Range("H6").Value = 0.438766
$ws.Range("I6").Value = 0.03816166876479336
$ws.Range("J6").Value = 0.03816166876479336
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4642756666666667
$ws.Range("N6").Value = 1.392827
$ws.Range("O6").Value = 0.5486456655941866
$ws.Range("P6").Value = 0.5486456655941866
$ws.Range("Q6").Value = 0.06790279238688889
$ws.Range("R6").Value = 0.611125131482
$ws.Range("S6").Value = 0.02093723415964493
$ws.Range("T6").Value = 0.02093723415964493

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cck"
$ws.Range("C7").Value = "Cckar"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1462553333333333
$ws.Range("H7").Value = 0.438766
$ws.Range("I7").Value = 0.03816166876479336
$ws.Range("J7").Value = 0.03816166876479336
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3819456666666667
$ws.Range("N7").Value = 1.145837
$ws.Range("O7").Value = 0.4513543344058135
$ws.Range("P7").Value = 0.4513543344058135
$ws.Range("Q7").Value = 0.05586159079355555
$ws.Range("R7").Value = 0.502754317142
$ws.Range("S7").Value = 0.01722443460514843
$ws.Range("T7").Value = 0.01722443460514843

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cck"
$ws.Range("C8").Value = "Cckar"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.256539666666666
$ws.Range("H8").Value = 6.769619
$ws.Range("I8").Value = 0.5887875494952928
$ws.Range("J8").Value = 0.5887875494952928
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4642756666666667
$ws.Range("N8").Value = 1.392827
$ws.Range("O8").Value = 0.5486456655941866
$ws.Range("P8").Value = 0.5486456655941866
$ws.Range("Q8").Value = 1.047656458101444
$ws.Range("R8").Value = 9.428908122913
$ws.Range("S8").Value = 0.323035736986415
$ws.Range("T8").Value = 0.323035736986415

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cck"
$ws.Range("C9").Value = "Cckar"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.256539666666666
$ws.Range("H9").Value = 6.769619
$ws.Range("I9").Value = 0.5887875494952928
$ws.Range("J9").Value = 0.5887875494952928
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3819456666666667
$ws.Range("N9").Value = 1.145837
$ws.Range("O9").Value = 0.4513543344058135
$ws.Range("P9").Value = 0.4513543344058135
$ws.Range("Q9").Value = 0.8618755473447777
$ws.Range("R9").Value = 7.756879926102999
$ws.Range("S9").Value = 0.2657518125088778
$ws.Range("T9").Value = 0.2657518125088778
